# growth_binary/binary_growth_research.xlsx
# "new exp data and vis #6 #7"
# Add a new additive entry (Maltoheptaose) as row 11 on the "additives" sheet,
# leave a stray space marker in F29, and move the active selection to D19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("additives")

# New data row (matches the layout used by the existing rows 8/9: chemical,
# BiGG id, final conc [g/L], citation, stock, is-related/original concentration).
# Values are entered in the same order the shared-string table grew in the
# source workbook: chemical, BiGG id, original/stock concentration (G),
# citation (E), stock solution (F), then the numeric final concentration (D).
$ws.Range("A11").Value = "Maltoheptaose"
$ws.Range("B11").Value = "malthp"
$ws.Range("G11").Value = "500müM"
$ws.Range("E11").Value = "Shim,2009"
$ws.Range("F11").Value = "25 g/L"
$ws.Range("D11").Value = 0.57650000000000001

# Stray trailing cell further down the sheet.
$ws.Range("F29").Value = " "

# Move the sheet's active selection.
[void]$ws.Range("D19").Select()
